$d = $word.ActiveDocument

$d.Paragraphs.Item(1).Range.Text = '这几(ji3)天心里颇不宁(ning4)静。今晚在院子里坐着(zhe)乘(cheng2)凉(liang2)，忽然想起日日走过(guo4)的(de)荷塘，在这满月的(de)光里，总该另有一番(fan1)样子吧(ba1)。月亮渐渐地(de)升高了(le)，墙外马路上孩子们的(de)欢笑，已经听不见(jian4)了(le)；妻在屋里拍着(zhe)闰儿(er2)，迷迷糊(hu2)糊(hu2)地(de)哼着(zhe)眠歌。我悄悄地(de)披了(le)大(da4)衫，带上门出去。'
$d.Paragraphs.Item(2).Range.Text = '沿着(zhe)荷塘，是一条曲(qu3)折(she2)的(de)小煤屑路。这是一条幽僻的(de)路；白天也少人走，夜晚更(geng4)加寂寞。荷塘四面，长着(zhe)许多树，蓊[weng3]蓊[weng3]郁郁的(de)。路的(de)一旁，是些杨柳，和(huo4)一些不知道名字的(de)树。没(mo4)有月光的(de)晚上，这路上阴森森的(de)，有些怕人。今晚却很好，虽然月光也还(huan2)是淡淡的(de)。'
$d.Paragraphs.Item(3).Range.Text = '路上只(zhi1)我一个(ge3)人，背着(zhuo2)手踱着(zhe)。这一片(pian4)天地(de)好像是我的(de)；我也像超出了(le)平常的(de)自己，到了(le)另一个(ge3)世界里。我爱热闹，也爱冷静；爱群居，也爱独处(chu3)。'
$d.Paragraphs.Item(4).Range.Text = '像今晚上，一个(ge3)人在这苍茫的(de)月下，什(shen2)么都(du1)可(ke4)以想，什(shen2)么都(du1)可(ke4)以不想，便(bian4)觉(jue2)是个(ge3)自由的(de)人。白天里一定要做的(de)事，一定要说(shuo1)的(de)话，现 在都(du1)可(ke4)不理。这是独处(chu3)的(de)妙处(chu4)，我且受用这无边的(de)荷香月色(se4)好了(le)。'
$d.Paragraphs.Item(5).Range.Text = '曲(qu3)曲(qu1)折(she2)折(she2)的(de)荷塘上面，弥望的(de)是田田的(de)叶子。叶子出水很高，像亭亭的(de)舞女的(de)裙。层层的(de)叶子中(zhong1)间(jian4)，零星地(de)点缀着(zhe)些白花，有袅[niao3]娜(nuo2)地(de)开着(zhe)的(de)，有羞涩[se4]地(de)打着(zhe)朵儿(er2)的(de)；正(zheng1)如一粒粒的(de)明珠，又如碧天里的(de)星星，又如刚出浴的(de)美人。微风过(guo4)处(chu3)，送来缕缕清香，仿佛(fu2)远处(chu4)高楼上渺[miao3]茫的(de)歌声似(si4)的(de)。'
$d.Paragraphs.Item(6).Range.Text = '忽然想起采莲的(de)事情来了(le)。采莲是江南的(de)旧俗，似(si4)乎很早就有，而六朝(chao2)时为盛；从诗歌里可(ke4)以约略知道。采莲的(de)是少年的(de)女子，她们是荡着(zhe)小船，唱着(zhe)艳歌去的(de)。采莲人不用说(shuo1)很多，还(huan2)有看(kan4)采莲的(de)人。那是一个(ge3)热闹的(de)季节(jie2)，也是一个(ge3)风流的(de)季节(jie2)。梁元帝《采莲赋》里说(shuo1)得(de)好：'
$d.Paragraphs.Item(7).Range.Text = '于是妖童媛[yuan2](yuan4)女，荡舟心许；鷁首徐回，兼传(chuan2)羽杯；棹[zhuo1]将(jiang4)移而藻挂，船欲动而萍开。尔其纤(qian4)腰束素，迁延顾步；夏始春余，叶嫩花初，恐沾裳而浅(qian3)笑，畏倾船而敛裾[ju1]。'
$d.Paragraphs.Item(8).Range.Text = '可(ke3)见(xian4)当(dang4)时嬉[xi1]游的(de)光景了(le)。这真是有趣的(de)事，可(ke3)惜我们现在早已无福消受了(le)。'
$d.Paragraphs.Item(10).Range.Text = '他勤奋好学，总是取得(de)好成绩。'
$d.Paragraphs.Item(11).Range.Text = '他怀着(zhe)沉重(chong2)的(de)心情，重(zhong4)新出发(fa4)。'
$d.Paragraphs.Item(12).Range.Text = '干(gan1)一行(xing2)爱一行(xing2)，你一定行(xing2)的(de)。'
$d.Paragraphs.Item(13).Range.Text = '知道了(le)，我了(liao3)解(jie3)这个(ge3)事情的(de)严重(chong2)性。'
